$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update J2/J3 text values (net effect of the shared-string reshuffle in the diff):
# J2 was "ExcelReport-NewImportLogic_3 - Test_Automation_3-Clinical-"
#   -> becomes "ExcelReport-NewImportLogic_3-Test_Automation_3-Clinical-"
# J3 was "WordReport-NewImportLogic_3 - Test_Automation_3-Clinical-" (unchanged text)
$ws.Range("J2").Value = "ExcelReport-NewImportLogic_3-Test_Automation_3-Clinical-"
$ws.Range("J3").Value = "WordReport-NewImportLogic_3 - Test_Automation_3-Clinical-"

# Update the view: scroll so H1 is the top-left cell, and select J4
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("J4").Select()

# Adjust window size (windowHeight 12720 -> 12576)
$excel.ActiveWindow.Height = $excel.ActiveWindow.Height - 144
